$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 40.75339133333333
$ws.Range("H2").Value = 122.260174
$ws.Range("I2").Value = 0.02126536631186857
$ws.Range("J2").Value = 0.02126536631186857
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.99315966666667
$ws.Range("N2").Value = 32.979479
$ws.Range("O2").Value = 0.5822520966482423
$ws.Range("P2").Value = 0.5822520966482423
$ws.Range("Q2").Value = 448.0085378854828
$ws.Range("R2").Value = 4032.076840969346
$ws.Range("S2").Value = 0.01238180412107837
$ws.Range("T2").Value = 0.01238180412107837
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 40.75339133333333
$ws.Range("H3").Value = 122.260174
$ws.Range("I3").Value = 0.02126536631186857
$ws.Range("J3").Value = 0.02126536631186857
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.231633
$ws.Range("N3").Value = 9.694898999999999
$ws.Range("O3").Value = 0.1711632639661454
$ws.Range("P3").Value = 0.1711632639661454
$ws.Range("Q3").Value = 131.700004294714
$ws.Range("R3").Value = 1185.300038652426
$ws.Range("S3").Value = 0.003639849507375135
$ws.Range("T3").Value = 0.003639849507375135
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 40.75339133333333
$ws.Range("H4").Value = 122.260174
$ws.Range("I4").Value = 0.02126536631186857
$ws.Range("J4").Value = 0.02126536631186857
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.9576306666666666
$ws.Range("N4").Value = 2.872892
$ws.Range("O4").Value = 0.05072085554911168
$ws.Range("P4").Value = 0.05072085554911168
$ws.Range("Q4").Value = 39.02669731146755
$ws.Range("R4").Value = 351.240275803208
$ws.Range("S4").Value = 0.001078597572903231
$ws.Range("T4").Value = 0.001078597572903231
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 40.75339133333333
$ws.Range("H5").Value = 122.260174
$ws.Range("I5").Value = 0.02126536631186857
$ws.Range("J5").Value = 0.02126536631186857
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.697989
$ws.Range("N5").Value = 11.093967
$ws.Range("O5").Value = 0.1958637838365006
$ws.Range("P5").Value = 0.1958637838365006
$ws.Range("Q5").Value = 150.705592863362
$ws.Range("R5").Value = 1356.350335770258
$ws.Range("S5").Value = 0.004165115110511827
$ws.Range("T5").Value = 0.004165115110511827
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1689.289306666667
$ws.Range("H6").Value = 5067.86792
$ws.Range("I6").Value = 0.8814813868902838
$ws.Range("J6").Value = 0.8814813868902838
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.99315966666667
$ws.Range("N6").Value = 32.979479
$ws.Range("O6").Value = 0.5822520966482423
$ws.Range("P6").Value = 0.5822520966482423
$ws.Range("Q6").Value = 18570.62707137929
$ws.Range("R6").Value = 167135.6436424137
$ws.Range("S6").Value = 0.5132443856732681
$ws.Range("T6").Value = 0.5132443856732681
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1689.289306666667
$ws.Range("H7").Value = 5067.86792
$ws.Range("I7").Value = 0.8814813868902838
$ws.Range("J7").Value = 0.8814813868902838
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.231633
$ws.Range("N7").Value = 9.694898999999999
$ws.Range("O7").Value = 0.1711632639661454
$ws.Range("P7").Value = 0.1711632639661454
$ws.Range("Q7").Value = 5459.16306997112
$ws.Range("R7").Value = 49132.46762974007
$ws.Range("S7").Value = 0.1508772313055456
$ws.Range("T7").Value = 0.1508772313055456
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1689.289306666667
$ws.Range("H8").Value = 5067.86792
$ws.Range("I8").Value = 0.8814813868902838
$ws.Range("J8").Value = 0.8814813868902838
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.9576306666666666
$ws.Range("N8").Value = 2.872892
$ws.Range("O8").Value = 0.05072085554911168
$ws.Range("P8").Value = 0.05072085554911168
$ws.Range("Q8").Value = 1617.715244936071
$ws.Range("R8").Value = 14559.43720442464
$ws.Range("S8").Value = 0.04470949009369271
$ws.Range("T8").Value = 0.04470949009369271
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1689.289306666667
$ws.Range("H9").Value = 5067.86792
$ws.Range("I9").Value = 0.8814813868902838
$ws.Range("J9").Value = 0.8814813868902838
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.697989
$ws.Range("N9").Value = 11.093967
$ws.Range("O9").Value = 0.1958637838365006
$ws.Range("P9").Value = 0.1958637838365006
$ws.Range("Q9").Value = 6246.973273870959
$ws.Range("R9").Value = 56222.75946483863
$ws.Range("S9").Value = 0.1726502798177773
$ws.Range("T9").Value = 0.1726502798177773
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 100.9654023333333
$ws.Range("H10").Value = 302.896207
$ws.Range("I10").Value = 0.05268435816499466
$ws.Range("J10").Value = 0.05268435816499466
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 10.99315966666667
$ws.Range("N10").Value = 32.979479
$ws.Range("O10").Value = 0.5822520966482423
$ws.Range("P10").Value = 0.5822520966482423
$ws.Range("Q10").Value = 1109.928788659572
$ws.Range("R10").Value = 9989.359097936152
$ws.Range("S10").Value = 0.03067557800213508
$ws.Range("T10").Value = 0.03067557800213508
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 100.9654023333333
$ws.Range("H11").Value = 302.896207
$ws.Range("I11").Value = 0.05268435816499466
$ws.Range("J11").Value = 0.05268435816499466
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.231633
$ws.Range("N11").Value = 9.694898999999999
$ws.Range("O11").Value = 0.1711632639661454
$ws.Range("P11").Value = 0.1711632639661454
$ws.Range("Q11").Value = 326.283126038677
$ws.Range("R11").Value = 2936.548134348093
$ws.Range("S11").Value = 0.009017626703481926
$ws.Range("T11").Value = 0.009017626703481926
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 100.9654023333333
$ws.Range("H12").Value = 302.896207
$ws.Range("I12").Value = 0.05268435816499466
$ws.Range("J12").Value = 0.05268435816499466
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.9576306666666666
$ws.Range("N12").Value = 2.872892
$ws.Range("O12").Value = 0.05072085554911168
$ws.Range("P12").Value = 0.05072085554911168
$ws.Range("Q12").Value = 96.68756554673821
$ws.Range("R12").Value = 870.188089920644
$ws.Range("S12").Value = 0.002672195720184357
$ws.Range("T12").Value = 0.002672195720184357
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 100.9654023333333
$ws.Range("H13").Value = 302.896207
$ws.Range("I13").Value = 0.05268435816499466
$ws.Range("J13").Value = 0.05268435816499466
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.697989
$ws.Range("N13").Value = 11.093967
$ws.Range("O13").Value = 0.1958637838365006
$ws.Range("P13").Value = 0.1958637838365006
$ws.Range("Q13").Value = 373.368947209241
$ws.Range("R13").Value = 3360.320524883169
$ws.Range("S13").Value = 0.01031895773919329
$ws.Range("T13").Value = 0.01031895773919329
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 85.41274733333334
$ws.Range("H14").Value = 256.238242
$ws.Range("I14").Value = 0.04456888863285297
$ws.Range("J14").Value = 0.04456888863285297
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 10.99315966666667
$ws.Range("N14").Value = 32.979479
$ws.Range("O14").Value = 0.5822520966482423
$ws.Range("P14").Value = 0.5822520966482423
$ws.Range("Q14").Value = 938.9559690039908
$ws.Range("R14").Value = 8450.603721035917
$ws.Range("S14").Value = 0.02595032885176065
$ws.Range("T14").Value = 0.02595032885176065
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 85.41274733333334
$ws.Range("H15").Value = 256.238242
$ws.Range("I15").Value = 0.04456888863285297
$ws.Range("J15").Value = 0.04456888863285297
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.231633
$ws.Range("N15").Value = 9.694898999999999
$ws.Range("O15").Value = 0.1711632639661454
$ws.Range("P15").Value = 0.1711632639661454
$ws.Range("Q15").Value = 276.022652903062
$ws.Range("R15").Value = 2484.203876127558
$ws.Range("S15").Value = 0.007628556449742748
$ws.Range("T15").Value = 0.007628556449742748
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 85.41274733333334
$ws.Range("H16").Value = 256.238242
$ws.Range("I16").Value = 0.04456888863285297
$ws.Range("J16").Value = 0.04456888863285297
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.9576306666666666
$ws.Range("N16").Value = 2.872892
$ws.Range("O16").Value = 0.05072085554911168
$ws.Range("P16").Value = 0.05072085554911168
$ws.Range("Q16").Value = 81.79386617065155
$ws.Range("R16").Value = 736.144795535864
$ws.Range("S16").Value = 0.002260572162331381
$ws.Range("T16").Value = 0.002260572162331381
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 85.41274733333334
$ws.Range("H17").Value = 256.238242
$ws.Range("I17").Value = 0.04456888863285297
$ws.Range("J17").Value = 0.04456888863285297
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.697989
$ws.Range("N17").Value = 9.694898999999999
$ws.Range("O17").Value = 0.1958637838365006
$ws.Range("P17").Value = 0.1958637838365006
$ws.Range("Q17").Value = 315.855400098446
$ws.Range("R17").Value = 2842.698600886014
$ws.Range("S17").Value = 0.008729431169018182
$ws.Range("T17").Value = 0.008729431169018182
